# Update the "K" column (column G) values for rows 2-13 to reflect
# regenerated save_data using K instead of Strike# (std/mean recalculated).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 2
    4  = 2
    5  = 3
    6  = 1
    7  = 3
    8  = 3
    9  = 5
    10 = 4
    11 = 4
    12 = 4
    13 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
